$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("B3").Value = 5900
$ws.Range("C3").Value = 9900
$ws.Range("D3").Value = 17640
$ws.Range("E3").Value = 32910
$ws.Range("F3").Value = 61960
$ws.Range("G3").Value = 123790
$ws.Range("H3").Value = 233340
$ws.Range("I3").Value = 463230
$ws.Range("J3").Value = 916310

$ws.Range("B3:J3").Select()
